$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "99.162.69"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Value = "3.300.71"
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").Value = "'254.72"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.91%  "
$ws.Range("D6").Value = "'624.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.71%  "
$ws.Range("D7").Value = "'1.45"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +25.91%  "
$ws.Range("E8").Value = "  +5.99%  "
$ws.Range("D10").Value = "'0.981"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +23.40%  "
$ws.Range("D11").Value = "3.298.50"
$ws.Range("E11").Value = "  -1.35%  "
$ws.Range("E12").Value = "  +0.58%  "
$ws.Range("D13").Value = "'39.71"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +11.02%  "
$ws.Range("D14").Value = "98.721.63"
$ws.Range("E14").Value = "  +0.96%  "
$ws.Range("D15").Value = "'0.0000249"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.56%  "
$ws.Range("D16").Value = "3.928.30"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("E17").Value = "  -0.70%  "
$ws.Range("D18").Value = "3.304.92"
$ws.Range("E18").Value = "  -1.26%  "
$ws.Range("E19").Value = "  -4.72%  "
$ws.Range("D20").Value = "'15.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.33%  "
$ws.Range("D21").Value = "'6.36"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +9.15%  "
$ws.Range("D22").Value = "'489.04"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.84%  "
$ws.Range("E23").Value = "  +1.11%  "
$ws.Range("E24").Value = "  -1.20%  "
$ws.Range("D25").Value = "'5.65"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "
$ws.Range("D26").Value = "'89.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.30%  "
$ws.Range("B27").Value = "Stellar"
$ws.Range("C27").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D27").Value = "'0.322"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +32.91%  "
$ws.Range("B28").Value = "Aptos"
$ws.Range("C28").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D28").Value = "'12.02"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.51%  "
$ws.Range("D29").Value = "3.487.36"
$ws.Range("E29").Value = "  -2.73%  "
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "'0.138"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +13.28%  "
$ws.Range("D32").Value = "'0.190"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.24%  "
$ws.Range("D33").Value = "'10.31"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +11.28%  "
$ws.Range("D34").Value = "'0.997"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.24%  "
$ws.Range("D35").Value = "'27.97"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.65%  "
$ws.Range("D36").Value = "'0.475"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +5.83%  "
$ws.Range("D37").Value = "'0.149"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").Value = "'7.22"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.06%  "
$ws.Range("D39").Value = "'1.95"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.53%  "
$ws.Range("D40").Value = "'24.80"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.25%  "
$ws.Range("D41").Value = "'490.31"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.83%  "
$ws.Range("D42").Value = "'3.62"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.96%  "
$ws.Range("E43").Value = "  -3.08%  "
$ws.Range("D44").Value = "'0.784"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.22%  "
$ws.Range("E45").Value = "  +0.00%  "
$ws.Range("D46").Value = "'3.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.77%  "
$ws.Range("B47").Value = "Monero"
$ws.Range("C47").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D47").Value = "'159.24"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.99%  "
$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.77%  "
$ws.Range("E49").Value = "  +16.16%  "
$ws.Range("E50").Value = "  +5.51%  "
$ws.Range("D51").Value = "'4.74"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +5.36%  "
